$d = $word.ActiveDocument

# 1. Delete the paragraph "Cada empleado debe estar asociado a un puesto definido en la tabla Puestos."
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Cada empleado debe estar asociado a un puesto definido en la tabla Puestos*") {
        $p.Range.Delete()
        break
    }
}

# 2. "Solo los empleados con puesto veterinario pueden ser asignados a un turno médico."
#    -> "Los empleados y veterinarios pueden asignar turnos médico."
$d.Content.Find.Execute(
    "Solo los empleados con puesto veterinario pueden ser asignados a un turno médico.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Los empleados y veterinarios pueden asignar turnos médico.", 2) | Out-Null

# 3. Insert two new paragraphs after "Toda entrada o salida de stock debe registrar la fecha y observaciones."
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Toda entrada o salida de stock debe registrar la fecha y observaciones*") {
        $p.Range.InsertParagraphAfter()
        $newPara = $p.Next()
        $newPara.Range.Text = "Cuando una sucursal registre una compra, automáticamente se debe actualizar el stock con los productos compras."
        $newPara.Range.InsertParagraphAfter()
        break
    }
}

# 4. "Toda venta debe generar un registro en Facturación con el detalle de los productos vendidos, cantidad y precio unitario."
#    -> "Toda venta debe generar un registro con el detalle de los productos vendidos, cantidad y precio unitario."
$d.Content.Find.Execute(
    "Toda venta debe generar un registro en Facturación con el detalle de los productos vendidos, cantidad y precio unitario.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Toda venta debe generar un registro con el detalle de los productos vendidos, cantidad y precio unitario.", 2) | Out-Null

# 5. "El total de la factura debe corresponder al detalle de los productos facturados."
#    -> "El total de la factura debe corresponder al detalle de los productos vendidos."
$d.Content.Find.Execute(
    "El total de la factura debe corresponder al detalle de los productos facturados.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El total de la factura debe corresponder al detalle de los productos vendidos.", 2) | Out-Null

# 6. "Una factura debe tener como mínimo un producto asociado."
#    -> "Una venta debe tener como mínimo un producto asociado."
$d.Content.Find.Execute(
    "Una factura debe tener como mínimo un producto asociado.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Una venta debe tener como mínimo un producto asociado.", 2) | Out-Null

# 7. Heading "Proveedores y Compras (Farmacias)" -> "Proveedores, Compras y Detalles compras"
$d.Content.Find.Execute(
    "Proveedores y Compras (Farmacias)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Proveedores, Compras y Detalles compras", 2) | Out-Null

# 8. "Toda compra realizada a una farmacia debe estar asociada a un empleado y a la farmacia correspondiente."
#    -> "Toda compra realizada en una sucursal debe estar asociada a un empleado "
$d.Content.Find.Execute(
    "Toda compra realizada a una farmacia debe estar asociada a un empleado y a la farmacia correspondiente.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Toda compra realizada en una sucursal debe estar asociada a un empleado ", 2) | Out-Null

# 9. Delete the paragraph "No se permite registrar compras con productos inexistentes en la tabla Productos."
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*No se permite registrar compras con productos inexistentes en la tabla Productos*") {
        $p.Range.Delete()
        break
    }
}

# 10. "Los precios de productos por farmacia deben actualizarse periódicamente en la tabla Farmacias_Productos."
#     -> "Los precios de productos deben actualizarse periódicamente en la tabla Productos"
$d.Content.Find.Execute(
    "Los precios de productos por farmacia deben actualizarse periódicamente en la tabla Farmacias_Productos.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Los precios de productos deben actualizarse periódicamente en la tabla Productos", 2) | Out-Null

# 11. Remove the stale "_GoBack" bookmark left over from the previous edit location.
try {
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks("_GoBack").Delete()
    }
} catch {
    Write-Host "bookmark cleanup skipped:" $_
}

Write-Host "Done"
